$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right marking 4 -> 5, Wrong marking -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Update "Total" row (row 12): recalc totals based on new marking scheme
$ws.Range("B12").Value = 95
$ws.Range("C12").Value = -7.199999999999999

# Update the Max score summary text in E12
$ws.Range("E12").Value = "87.8/140"
